$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data

# Row 2
$ws.Range("D2").Value = "67.778.69"
$ws.Range("E2").Value = "  +3.24%  "

# Row 3
$ws.Range("D3").Value = "3.283.09"
$ws.Range("E3").Value = "  +0.15%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.34%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.59"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.09%  "

# Row 6
$ws.Range("E6").Value = "  -1.21%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.19%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.584"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.94%  "

# Row 9
$ws.Range("D9").Value = "3.279.02"
$ws.Range("E9").Value = "  +0.49%  "

# Row 10
$ws.Range("E10").Value = "  +1.16%  "

# Row 11
$ws.Range("E11").Value = "  +1.39%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.70"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.51%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000270"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.68%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "691.72"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +13.44%  "

# Row 15
$ws.Range("D15").Value = "3.811.48"
$ws.Range("E15").Value = "  -0.32%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.37"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.06%  "

# Row 17
$ws.Range("D17").Value = "67.867.61"
$ws.Range("E17").Value = "  +3.41%  "

# Row 18
$ws.Range("E18").Value = "  +1.48%  "

# Row 19
$ws.Range("D19").Value = "3.285.64"
$ws.Range("E19").Value = "  -0.51%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.41"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.66%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.83"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.29%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.893"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.85%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.17"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.17"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.67%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.29"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.07%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.95"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.37%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.74"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.09%  "

# Row 28
$ws.Range("E28").Value = "  -2.65%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.40"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.20%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.77"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.82%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.44"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.46%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.73"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.98%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "586.44"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.70%  "

# Row 34
$ws.Range("D34").Value = "3.886.89"
$ws.Range("E34").Value = "  +2.54%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.86"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.69%  "

# Row 36
$ws.Range("E36").Value = "  +1.18%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.40"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -8.76%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.38"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.130"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.25%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.25"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.19%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.63"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.75%  "

# Row 43
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.41"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.44%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "32.31"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.05%  "

# Row 45
$ws.Range("D45").Value = "0.0₃0676"
$ws.Range("E45").Value = "  +0.45%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.331"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.55%  "

# Row 47
$ws.Range("E47").Value = "  +2.19%  "

# Row 48
$ws.Range("E48").Value = "  +2.05%  "

# Row 49
$ws.Range("E49").Value = "  +0.44%  "

# Row 50
$ws.Range("E50").Value = "  +10.50%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.50"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.56%  "
